$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value2 = 45043
$ws.Range("L2").Value2 = "Primera"
$ws.Range("N2").Value2 = 19000
$ws.Range("O2").Value2 = 20000
$ws.Range("P2").Value2 = 19500
$ws.Range("R2").Value2 = "Región de O'Higgins"
$ws.Range("S2").Value2 = 1083

# Row 4 updates
$ws.Range("D4").Value2 = 45086
$ws.Range("L4").Value2 = "Segunda"
$ws.Range("N4").Value2 = 20000
$ws.Range("O4").Value2 = 21000
$ws.Range("P4").Value2 = 20500
$ws.Range("R4").Value2 = "Provincia de Curicó"
$ws.Range("S4").Value2 = 1139
